# Applies the ElderScrollsExplorer ScrollsExplorer.docx edit described by the diff:
#  1. The first bullet ("Using obliv spiderdaedra ...") is replaced with a new
#     bullet about the .kf file memory usage.
#  2. Six new bullets (TES3 animations, BSparticle, water colors, water vertex
#     shader question, fade system question, and an empty "_GoBack" bookmark
#     paragraph) are inserted, followed by the ORIGINAL "Using obliv
#     spiderdaedra ..." bullet text (moved further down the list).
#  3. A new bullet ("Fo3 has black textures where should be transparent") is
#     inserted after the "...and it don't play nice" bullet.
#  4. Three pre-existing runs gain a <w:lastRenderedPageBreak/> marker.
#
# NOTE: paragraph/range COM handles returned by this host are snapshots tied
# to a position - they do NOT stay "live" across a later mutation elsewhere
# in the document. So every step below re-fetches $d.Paragraphs(<index>)
# immediately before using it, instead of caching paragraph objects across
# InsertParagraphBefore/InsertXML calls.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function New-ListParagraphXml([string]$innerXml) {
    return '<w:p ' + $wNs + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0"/></w:pPr>' + $innerXml + '</w:p>'
}

function Insert-EmptyListParagraphsBefore([int]$index, [int]$count) {
    # Inserts $count brand-new empty ListParagraph paragraphs immediately
    # before paragraph $index (1-based). After this, the original paragraph
    # that used to be at $index is now at $index + $count.
    for ($n = 0; $n -lt $count; $n++) {
        $d.Paragraphs($index).Range.InsertParagraphBefore()
    }
}

function Set-ParagraphXml([int]$index, [string]$innerXml) {
    $d.Paragraphs($index).Range.InsertXML((New-ListParagraphXml $innerXml))
}

# ============================================================================
# Part A: rework the first bullet + insert the six new bullets before it
# ============================================================================

# Make room for 7 new paragraphs before the "Using obliv ..." bullet
# (currently paragraph 2). After this, that original paragraph sits at 9.
Insert-EmptyListParagraphsBefore 2 7

$kfXml = '<w:r><w:t xml:space="preserve">The </w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/><w:r><w:t>kf</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
         '<w:r><w:t xml:space="preserve"> file keep loading in increasing memory use, perhaps I should auto load all of them when a CREA or NPC_ is loaded to fix the memory usage (and unload them of course)</w:t></w:r>'
Set-ParagraphXml 2 $kfXml

$tes3Anim = '<w:r><w:t xml:space="preserve">TES3: </w:t></w:r>' +
            '<w:r><w:t xml:space="preserve">Animations on </w:t></w:r>' +
            '<w:proofErr w:type="spellStart"/><w:r><w:t>crea</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
            '<w:r><w:t xml:space="preserve"> and </w:t></w:r>' +
            '<w:proofErr w:type="spellStart"/><w:r><w:t>npc</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
            '<w:r><w:t xml:space="preserve"> not going tes3</w:t></w:r>'
Set-ParagraphXml 3 $tes3Anim

$bsParticle = '<w:r><w:t xml:space="preserve">TES3: </w:t></w:r>' +
              '<w:proofErr w:type="spellStart"/><w:r><w:t>BSparticle</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
              '<w:r><w:t xml:space="preserve"> system not working</w:t></w:r>'
Set-ParagraphXml 4 $bsParticle

$waterColors = '<w:r><w:t xml:space="preserve">Water not overriding </w:t></w:r>' +
               '<w:proofErr w:type="spellStart"/><w:r><w:t>colors</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
               '<w:r><w:t xml:space="preserve"> for </w:t></w:r>' +
               '<w:proofErr w:type="spellStart"/><w:r><w:t>skyrim</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
               '<w:r><w:t xml:space="preserve"> and fallout</w:t></w:r>'
Set-ParagraphXml 5 $waterColors

$waterVert = '<w:r><w:t>Could water use the 0th row as the last row for vertex shaking to make perfect match ups?</w:t></w:r>' +
             '<w:r><w:t xml:space="preserve"> In the </w:t></w:r>' +
             '<w:proofErr w:type="spellStart"/><w:r><w:t>water.vert</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
             '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
             '<w:proofErr w:type="spellStart"/><w:r><w:t>shader</w:t></w:r><w:proofErr w:type="spellEnd"/>'
Set-ParagraphXml 6 $waterVert

$fadeSystem = '<w:r><w:t>For the fade system surely telling the scene graph to skip the BG is better than detaching and re attaching (expensive surely)</w:t></w:r>'
Set-ParagraphXml 7 $fadeSystem

$goBack = '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
Set-ParagraphXml 8 $goBack

$origRestored = '<w:r><w:t xml:space="preserve">Using </w:t></w:r>' +
                '<w:proofErr w:type="spellStart"/><w:r><w:t>obliv</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
                '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
                '<w:proofErr w:type="spellStart"/><w:r><w:t>spiderdaedra</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
                '<w:r><w:t xml:space="preserve"> I see skin failing to map quite badly, I presume that bones </w:t></w:r>' +
                '<w:r><w:t>fail</w:t></w:r>' +
                '<w:r><w:t xml:space="preserve"> too </w:t></w:r>' +
                '<w:r><w:t>in a different bug</w:t></w:r>' +
                '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
                '<w:r><w:t xml:space="preserve">(tight shoulder thing). This skin failure is also shown in the </w:t></w:r>' +
                '<w:proofErr w:type="spellStart"/><w:r><w:t>skyrim</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
                '<w:r><w:t xml:space="preserve"> cock ups</w:t></w:r>'
Set-ParagraphXml 9 $origRestored

# At this point paragraph 10 is "JInternalFramse for DisplayDialog..." (the
# original paragraph 3), unchanged, and paragraph 11 is "The physics
# synchronized calls now make the KCC pause on load of new land!" (original
# paragraph 4).

# ============================================================================
# Part B: insert the "Fo3 has black textures ..." bullet right after the
#         "...and it don't play nice" bullet (paragraph 10), i.e. before the
#         "physics synchronized" bullet (paragraph 11).
# ============================================================================

Insert-EmptyListParagraphsBefore 11 1
$fo3Black = '<w:r><w:t>Fo3 has black textures where should be transparent</w:t></w:r>'
Set-ParagraphXml 11 $fo3Black

Write-Host "part-b-done"

# ============================================================================
# Part C: add <w:lastRenderedPageBreak/> to the first run of three existing
#         paragraphs (their content is otherwise untouched - the exact
#         original markup, including rsid/proofErr detail, is reproduced so
#         only the page-break marker is new).
# ============================================================================

function Find-ParagraphIndexContaining([string]$needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs($i).Range.Text -like "*$needle*") {
            return $i
        }
    }
    return -1
}

$glslParaXml = '<w:p w:rsidR="00966725" w:rsidRDefault="00966725" w:rsidP="00966725"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">Still haven' + [char]0x2019 + 't got GLSL FBO </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>shader</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> working, must use the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>offscreen</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> code but make a front buffer and shade and flip like web examples do</w:t></w:r></w:p>'
$glslIdx = Find-ParagraphIndexContaining "Still haven"
$d.Paragraphs($glslIdx).Range.InsertXML($glslParaXml)

$launcherParaXml = '<w:p w:rsidR="00F70058" w:rsidRPr="00F70058" w:rsidRDefault="00F70058" w:rsidP="00D553D6"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Times New Roman"/><w:lang w:eastAsia="en-NZ"/></w:rPr></w:pPr><w:r w:rsidRPr="00F70058"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Times New Roman"/><w:lang w:eastAsia="en-NZ"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">Launcher could then include a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00F70058"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Times New Roman"/><w:lang w:eastAsia="en-NZ"/></w:rPr><w:t>url</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00F70058"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Times New Roman"/><w:lang w:eastAsia="en-NZ"/></w:rPr><w:t xml:space="preserve"> to get latest from and unzip over the current jar, this probably means moving the launcher into a separate jar file to not over ride itself.</w:t></w:r></w:p>'
$launcherIdx = Find-ParagraphIndexContaining "Launcher could then include a"
$d.Paragraphs($launcherIdx).Range.InsertXML($launcherParaXml)

$jvmParaXml = '<w:p w:rsidR="008256A3" w:rsidRDefault="008256A3" w:rsidP="008256A3"><w:pPr><w:pStyle w:val="NormalWeb"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>The first solution you could choose would be to install the files in the extensions folder of the JVM (JAVA_HOME\lib\</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ext</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>). The problem with this approach is that your application will only work until a new version of Java is installed.</w:t></w:r></w:p>'
$jvmIdx = Find-ParagraphIndexContaining "The first solution you could choose"
$d.Paragraphs($jvmIdx).Range.InsertXML($jvmParaXml)

Write-Host "part-c-done"
